$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("checklist")

# Add missing year values for rows 14 and 23
$ws.Range("A14").Value = 2007
$ws.Range("A23").Value = 2011

# Update the active cell selection on the sheet view
$ws.Activate()
$ws.Range("A24").Select()
